$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.631.41"
$ws.Range("E2").Value = "  -2.16%  "
$ws.Range("D3").Value = "'1.841.56"
$ws.Range("E3").Value = "  -1.39%  "
$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  -0.21%  "
$ws.Range("E5").Value = "  -1.41%  "
$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = "  -0.22%  "
$ws.Range("D7").Value = "'0.4261"
$ws.Range("E7").Value = "  -2.66%  "
$ws.Range("D8").Value = "'0.3645"
$ws.Range("E8").Value = "  -1.51%  "
$ws.Range("D9").Value = "'45.71"
$ws.Range("E9").Value = "  +1.30%  "
$ws.Range("D10").Value = "'0.07284"
$ws.Range("E10").Value = "  -3.00%  "
$ws.Range("D11").Value = "'0.8958"
$ws.Range("E11").Value = "  -4.52%  "
$ws.Range("D12").Value = "'20.60"
$ws.Range("E12").Value = "  -3.89%  "
$ws.Range("D13").Value = "'1.855.50"
$ws.Range("E13").Value = "  -3.21%  "
$ws.Range("D14").Value = "'5.380"
$ws.Range("E14").Value = "  -1.03%  "
$ws.Range("D15").Value = "'6.556"
$ws.Range("E15").Value = "  -2.44%  "
$ws.Range("D16").Value = "'0.06859"
$ws.Range("E16").Value = "  -0.18%  "
$ws.Range("D17").Value = "'1.002"
$ws.Range("E17").Value = "  -0.15%  "
$ws.Range("D18").Value = "'78.34"
$ws.Range("E18").Value = "  -4.60%  "
$ws.Range("D19").Value = "'0.000008829"
$ws.Range("E19").Value = "  -2.69%  "
$ws.Range("D20").Value = "'0.9993"
$ws.Range("E20").Value = "  -0.40%  "
$ws.Range("D21").Value = "'15.56"
$ws.Range("E21").Value = "  -2.08%  "
$ws.Range("D22").Value = "'27.640.84"
$ws.Range("E22").Value = "  -2.13%  "
$ws.Range("D23").Value = "'4.973"
$ws.Range("E23").Value = "  -3.00%  "
$ws.Range("D24").Value = "'10.58"
$ws.Range("E24").Value = "  -1.61%  "
$ws.Range("D25").Value = "'2.078.18"
$ws.Range("E25").Value = "  -3.37%  "
$ws.Range("E26").Value = "  +1.05%  "
$ws.Range("D27").Value = "'154.49"
$ws.Range("E27").Value = "  -0.19%  "
$ws.Range("D28").Value = "'18.26"
$ws.Range("E28").Value = "  -0.77%  "
$ws.Range("D29").Value = "'5.241"
$ws.Range("E29").Value = "  -1.24%  "
$ws.Range("D30").Value = "'1.832"
$ws.Range("E30").Value = "  +5.96%  "
$ws.Range("D31").Value = "'111.43"
$ws.Range("E31").Value = "  -2.22%  "
$ws.Range("D32").Value = "'0.08889"
$ws.Range("E32").Value = "  -1.67%  "
$ws.Range("D33").Value = "'0.7747"
$ws.Range("E33").Value = "  -2.93%  "
$ws.Range("B34").Value = "HuobiToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D34").Value = "'2.997"
$ws.Range("E34").Value = "  +1.39%  "
$ws.Range("B35").Value = "Filecoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D35").Value = "'4.567"
$ws.Range("E35").Value = "  -5.50%  "
$ws.Range("E36").Value = "  -6.21%  "
$ws.Range("D38").Value = "'0.05438"
$ws.Range("E38").Value = "  +0.20%  "
$ws.Range("D39").Value = "'1.098"
$ws.Range("E39").Value = "  -1.96%  "
$ws.Range("D40").Value = "'0.01926"
$ws.Range("E40").Value = "  -1.28%  "
$ws.Range("D41").Value = "'2.771"
$ws.Range("E41").Value = "  -5.45%  "
$ws.Range("D42").Value = "'0.5065"
$ws.Range("E42").Value = "  -3.56%  "
$ws.Range("D43").Value = "'6.800"
$ws.Range("E43").Value = "  -4.21%  "
$ws.Range("D44").Value = "'0.1644"
$ws.Range("E44").Value = "  -1.91%  "
$ws.Range("D45").Value = "'8.220"
$ws.Range("E45").Value = "  -5.66%  "
$ws.Range("D46").Value = "'0.06641"
$ws.Range("E46").Value = "  -1.62%  "
$ws.Range("D47").Value = "'10.42"
$ws.Range("E47").Value = "  -1.43%  "
$ws.Range("D48").Value = "'105.85"
$ws.Range("E48").Value = "  -1.84%  "
$ws.Range("D49").Value = "'0.4712"
$ws.Range("E49").Value = "  -3.24%  "
$ws.Range("D50").Value = "'1.000"
$ws.Range("E50").Value = "  -0.20%  "
$ws.Range("D51").Value = "'1.636"
$ws.Range("E51").Value = "  -2.45%  "
